$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 2
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 7
$ws.Range("M2").Value = 7

# Row 3
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 8
$ws.Range("K3").Value = 8
$ws.Range("L3").Value = 8
$ws.Range("M3").Value = 7

# Row 4 (new)
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 5
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

# Row 7
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 11
$ws.Range("G7").Value = 11
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 15
$ws.Range("L7").Value = 16
$ws.Range("M7").Value = 16

# Row 8
$ws.Range("C8").Value = 13
$ws.Range("D8").Value = 12
$ws.Range("E8").Value = 12
$ws.Range("F8").Value = 12
$ws.Range("G8").Value = 11
$ws.Range("I8:J8").ClearContents()
$ws.Range("K8").Value = 17
$ws.Range("L8").Value = 17
$ws.Range("M8").Value = 17

# Row 9 (new)
$ws.Range("C9").Value = 13
$ws.Range("D9").Value = 13
$ws.Range("E9").Value = 14
$ws.Range("F9").Value = 14
$ws.Range("G9").Value = 14

# Update selection to match target (G12)
$ws.Range("G12").Select()
